$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the two combo/ruby_eval output formulas so they coerce operands to
# integers before adding them together (all output fields used in the
# combo_order_eval table must have values).
$ws.Range("E5").Value = '${r:${BioTime price}.to_i + ${ProSpec price}.to_i}'
$ws.Range("D5").Value = '${r:${BioTime Handling}.to_i + ${ProSpec Handling}.to_i}'

# Widen column D slightly to fit the updated formula text.
$ws.Columns.Item(4).ColumnWidth = 56.8333333333

# Move the active selection to D5, matching the saved view state.
$ws.Range("D5").Select()
